$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the team lists in column O ("tied_teams") for several match rows.
$ws.Range("O3").Value = '[''Uruguay'', ''Hungary'']'
$ws.Range("O4").Value = '[''Uruguay'', ''Hungary'']'
$ws.Range("O9").Value = '[''Northern Ireland'', ''Bulgaria'']'
$ws.Range("O10").Value = '[''Northern Ireland'', ''Bulgaria'']'
$ws.Range("O11").Value = '[''Northern Ireland'', ''Bulgaria'']'
$ws.Range("O12").Value = '[''Northern Ireland'', ''Bulgaria'']'
$ws.Range("O40").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O41").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O42").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O43").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O44").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O45").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O46").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O47").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O48").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O49").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O50").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O51").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O52").Value = '[''Costa Rica'', ''Ireland'']'
$ws.Range("O53").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O54").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O55").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O56").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O57").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O58").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O59").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O60").Value = '[''Colombia'', ''Costa Rica'', ''Ireland'', ''Argentina'']'
$ws.Range("O61").Value = '[''Colombia'', ''Argentina'']'
$ws.Range("O62").Value = '[''Colombia'', ''Argentina'']'
$ws.Range("O63").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O64").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O65").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O66").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O67").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O68").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O69").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O70").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O71").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O72").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Range("O73").Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'